# Auto-generated edit script: updates recalculated market-price/profit
# figures across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (41 cell changes) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4168120
$ws.Range("J19").Value = 6668275.5
$ws.Range("L19").Value = 6668275.5
$ws.Range("N19").Value = -6668625.5
$ws.Range("H43").Value = 4256.4287
$ws.Range("I43").Value = 4256.4287
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 4256.4287
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -4187.4287
$ws.Range("N43").ClearContents()
$ws.Range("H86").Value = 7219.5835
$ws.Range("I86").Value = 3673.125
$ws.Range("K86").Value = 3673.125
$ws.Range("M86").Value = -2550.125
$ws.Range("H89").Value = 7219.5835
$ws.Range("I89").Value = 3673.125
$ws.Range("K89").Value = 18365.625
$ws.Range("M89").Value = -12749.625
$ws.Range("H92").Value = 3257.1
$ws.Range("I92").Value = 1907.1428
$ws.Range("J92").Value = 6407
$ws.Range("K92").Value = 1907.1428
$ws.Range("L92").Value = 6407
$ws.Range("M92").Value = -659.1428000000001
$ws.Range("N92").Value = -8903
$ws.Range("H101").Value = 2896.75
$ws.Range("I101").Value = 196
$ws.Range("K101").Value = 588
$ws.Range("M101").Value = 1034
$ws.Range("H107").Value = 3141.8572
$ws.Range("I107").Value = 666.3333
$ws.Range("K107").Value = 666.3333
$ws.Range("M107").Value = 1253.6667
$ws.Range("H116").Value = 11738.73
$ws.Range("I116").Value = 6850.846
$ws.Range("J116").Value = 16626.615
$ws.Range("K116").Value = 6850.846
$ws.Range("L116").Value = 16626.615
$ws.Range("M116").Value = -3408.846
$ws.Range("N116").Value = -23510.615

# --- Sheet: ARM (22 cell changes) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1558.1154
$ws.Range("I97").Value = 1707.7778
$ws.Range("J97").Value = 1221.375
$ws.Range("K97").Value = 1707.7778
$ws.Range("L97").Value = 1221.375
$ws.Range("M97").Value = -1211.7778
$ws.Range("N97").Value = -2213.375
$ws.Range("H114").Value = 119999.5
$ws.Range("J114").Value = 119999.5
$ws.Range("L114").Value = 119999.5
$ws.Range("N114").Value = -128677.5
$ws.Range("H122").Value = 3128
$ws.Range("I122").Value = 2840.7273
$ws.Range("K122").Value = 8522.1819
$ws.Range("M122").Value = -6072.1819
$ws.Range("H132").Value = 2132740
$ws.Range("I132").Value = 4177.1665
$ws.Range("J132").Value = 9098946
$ws.Range("K132").Value = 12531.4995
$ws.Range("L132").Value = 27296838
$ws.Range("M132").Value = -10001.4995
$ws.Range("N132").Value = -27301898

# --- Sheet: BSM (31 cell changes) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1067.2778
$ws.Range("J20").Value = 1247
$ws.Range("L20").Value = 1247
$ws.Range("N20").Value = -1741
$ws.Range("H22").Value = 940.4
$ws.Range("I22").Value = 1055.9231
$ws.Range("J22").Value = 189.5
$ws.Range("K22").Value = 1055.9231
$ws.Range("L22").Value = 189.5
$ws.Range("M22").Value = -882.9231
$ws.Range("N22").Value = -535.5
$ws.Range("H86").Value = 4969.5
$ws.Range("J86").Value = 7002
$ws.Range("L86").Value = 7002
$ws.Range("N86").Value = -9248
$ws.Range("H89").Value = 4969.5
$ws.Range("J89").Value = 7002
$ws.Range("L89").Value = 35010
$ws.Range("N89").Value = -46242
$ws.Range("H94").Value = 2017.8096
$ws.Range("I94").Value = 2073.6875
$ws.Range("K94").Value = 2073.6875
$ws.Range("M94").Value = -1622.6875
$ws.Range("H105").Value = 563701.9399999999
$ws.Range("I105").Value = 717769.8
$ws.Range("K105").Value = 717769.8
$ws.Range("M105").Value = -716022.8
$ws.Range("H128").Value = 6670.8335
$ws.Range("I128").Value = 6670.8335
$ws.Range("K128").Value = 20012.5005
$ws.Range("M128").Value = -17522.5005

# --- Sheet: CRP (26 cell changes) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26319318
$ws.Range("J31").Value = 3798.3125
$ws.Range("L31").Value = 3798.3125
$ws.Range("N31").Value = -4388.3125
$ws.Range("H34").Value = 26319318
$ws.Range("J34").Value = 3798.3125
$ws.Range("L34").Value = 3798.3125
$ws.Range("N34").Value = -4202.3125
$ws.Range("I107").Value = 483.25
$ws.Range("J107").Value = 1714.875
$ws.Range("K107").Value = 483.25
$ws.Range("L107").Value = 1714.875
$ws.Range("M107").Value = 1436.75
$ws.Range("N107").Value = -5554.875
$ws.Range("H122").Value = 3524.8235
$ws.Range("I122").Value = 3086.3845
$ws.Range("K122").Value = 9259.1535
$ws.Range("M122").Value = -6809.1535
$ws.Range("H132").Value = 3163.4092
$ws.Range("I132").Value = 1926.8948
$ws.Range("K132").Value = 5780.6844
$ws.Range("M132").Value = -3250.6844
$ws.Range("H134").Value = 1673.8889
$ws.Range("I134").Value = 1612.2
$ws.Range("K134").Value = 4836.6
$ws.Range("M134").Value = -2301.6

# --- Sheet: CUL (44 cell changes) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 128
$ws.Range("J2").Value = 65.833336
$ws.Range("L2").Value = 395.000016
$ws.Range("N2").Value = -621.000016
$ws.Range("H5").Value = 326.8125
$ws.Range("I5").Value = 265.54544
$ws.Range("J5").Value = 461.6
$ws.Range("K5").Value = 796.63632
$ws.Range("L5").Value = 1384.8
$ws.Range("M5").Value = -684.63632
$ws.Range("N5").Value = -1608.8
$ws.Range("H33").Value = 5774807.5
$ws.Range("I33").Value = 146.66667
$ws.Range("J33").Value = 7699694
$ws.Range("K33").Value = 880.0000200000001
$ws.Range("L33").Value = 46198164
$ws.Range("M33").Value = -597.0000200000001
$ws.Range("N33").Value = -46198730
$ws.Range("H98").Value = 717
$ws.Range("J98").Value = 641.5714
$ws.Range("L98").Value = 1924.7142
$ws.Range("N98").Value = -4920.7142
$ws.Range("H104").Value = 24166.5
$ws.Range("J104").Value = 24166.5
$ws.Range("L104").Value = 72499.5
$ws.Range("N104").Value = -77741.5
$ws.Range("H107").Value = 5358422.5
$ws.Range("J107").Value = 8279843
$ws.Range("L107").Value = 24839529
$ws.Range("N107").Value = -24843369
$ws.Range("H113").Value = 1437.3334
$ws.Range("I113").Value = 1465.2222
$ws.Range("J113").Value = 1409.4445
$ws.Range("K113").Value = 4395.6666
$ws.Range("L113").Value = 4228.333500000001
$ws.Range("M113").Value = -2225.6666
$ws.Range("N113").Value = -8568.333500000001
$ws.Range("H135").Value = 326.8125
$ws.Range("I135").Value = 265.54544
$ws.Range("J135").Value = 461.6
$ws.Range("K135").Value = 2389.90896
$ws.Range("L135").Value = 4154.400000000001
$ws.Range("M135").Value = 145.0910400000002
$ws.Range("N135").Value = -9224.400000000001

# --- Sheet: GSM (36 cell changes) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 26575.666
$ws.Range("J44").Value = 26575.666
$ws.Range("L44").Value = 26575.666
$ws.Range("N44").Value = -27767.666
$ws.Range("H70").Value = 9068.704
$ws.Range("I70").Value = 9501.5
$ws.Range("J70").Value = 8886.474
$ws.Range("K70").Value = 9501.5
$ws.Range("L70").Value = 8886.474
$ws.Range("M70").Value = -9231.5
$ws.Range("N70").Value = -9426.474
$ws.Range("H73").Value = 9068.704
$ws.Range("I73").Value = 9501.5
$ws.Range("J73").Value = 8886.474
$ws.Range("K73").Value = 9501.5
$ws.Range("L73").Value = 8886.474
$ws.Range("M73").Value = -8565.5
$ws.Range("N73").Value = -10758.474
$ws.Range("H97").Value = 1067.2778
$ws.Range("I97").Value = 1173.7273
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 1173.7273
$ws.Range("L97").Value = 900
$ws.Range("M97").Value = -677.7273
$ws.Range("N97").Value = -1892
$ws.Range("H107").Value = 394
$ws.Range("I107").Value = 365.2
$ws.Range("J107").Value = 414.57144
$ws.Range("K107").Value = 365.2
$ws.Range("L107").Value = 414.57144
$ws.Range("M107").Value = 1554.8
$ws.Range("N107").Value = -4254.57144
$ws.Range("H122").Value = 2834.3928
$ws.Range("I122").Value = 2613
$ws.Range("K122").Value = 7839
$ws.Range("M122").Value = -5389

# --- Sheet: LTW (33 cell changes) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 13599.556
$ws.Range("I22").Value = 23479.2
$ws.Range("K22").Value = 23479.2
$ws.Range("M22").Value = -23184.2
$ws.Range("H27").Value = 13599.556
$ws.Range("I27").Value = 23479.2
$ws.Range("K27").Value = 23479.2
$ws.Range("M27").Value = -23372.2
$ws.Range("H55").Value = 1018.2917
$ws.Range("I55").Value = 472.1
$ws.Range("J55").Value = 1408.4286
$ws.Range("K55").Value = 472.1
$ws.Range("L55").Value = 1408.4286
$ws.Range("M55").Value = -299.1
$ws.Range("N55").Value = -1754.4286
$ws.Range("H61").Value = 111115544
$ws.Range("I61").Value = 166670160
$ws.Range("J61").Value = 6297.6665
$ws.Range("K61").Value = 166670160
$ws.Range("L61").Value = 6297.6665
$ws.Range("M61").Value = -166669958
$ws.Range("N61").Value = -6701.6665
$ws.Range("H113").Value = 111115544
$ws.Range("I113").Value = 166670160
$ws.Range("J113").Value = 6297.6665
$ws.Range("K113").Value = 166670160
$ws.Range("L113").Value = 6297.6665
$ws.Range("M113").Value = -166667990
$ws.Range("N113").Value = -10637.6665
$ws.Range("H136").Value = 9998.666999999999
$ws.Range("I136").Value = 9998
$ws.Range("K136").Value = 29994
$ws.Range("M136").Value = -27444

# --- Sheet: WVR (25 cell changes) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 13335
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 13335
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 13335
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -15207
$ws.Range("H77").Value = 13335
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 13335
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 40005
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -49365
$ws.Range("H122").Value = 2753.24
$ws.Range("I122").Value = 2749.25
$ws.Range("K122").Value = 8247.75
$ws.Range("M122").Value = -5797.75
$ws.Range("H126").Value = 7869
$ws.Range("I126").Value = 9151
$ws.Range("J126").Value = 6074.2
$ws.Range("K126").Value = 27453
$ws.Range("L126").Value = 18222.6
$ws.Range("M126").Value = -24983
$ws.Range("N126").Value = -23162.6

Write-Host "Applied 258 cell updates across 8 sheets"